$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "One"
$ws.Range("D3").Value = "one"
$ws.Range("F6").Value = "one"

$ws.Range("F6").Select()
